$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.364.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4668'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06569'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07894'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.874.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.128'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6758'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '281.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.353.72'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.498'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.115.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007302'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.199'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.937'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.46%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.427'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.477'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.124'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04686'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.121'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7062'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.712'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01861'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.322'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.545'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.946'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8482'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.207'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.174'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '932.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("E51").Value = '  -2.46%  '
